$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 146-148)
$rows = @(
    @{ Row = 146; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44939; E = 8; F = 100112028; G = "Sandia"; H = "Sin especificar"; I = "Extra";    J = 500; K = 3000; L = 3000; M = 3000; N = "$/unidad"; O = "Región de O'Higgins"; P = 3000; Q = 1; R = "Hortaliza" },
    @{ Row = 147; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44939; E = 8; F = 100112028; G = "Sandia"; H = "Sin especificar"; I = "Primera";  J = 500; K = 2500; L = 2500; M = 2500; N = "$/unidad"; O = "Región de O'Higgins"; P = 2500; Q = 1; R = "Hortaliza" },
    @{ Row = 148; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44939; E = 8; F = 100112028; G = "Sandia"; H = "Sin especificar"; I = "Segunda"; J = 500; K = 2000; L = 2000; M = 2000; N = "$/unidad"; O = "Región de O'Higgins"; P = 2000; Q = 1; R = "Hortaliza" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C

    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
}
